$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 69: change C69's text from "Recursion" to "Leatcode"
$ws.Range("C69").Value = "Leatcode"

# New row 73
$ws.Range("A73").Value = "                  10/23/2024"
$ws.Range("B73").Value = "isSame "
$ws.Range("C73").Value = "Leatcode"
$ws.Range("G73").Value = "Easy"
$ws.Range("H73").Value = "Yes"

# Update selection to match diff
$ws.Range("J68").Select()
